$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 11465.556
$ws.Range("I51").Value = 11990
$ws.Range("J51").Value = 11400
$ws.Range("K51").Value = 11990
$ws.Range("L51").Value = 11400
$ws.Range("M51").Value = -11506
$ws.Range("N51").Value = -12368

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H61").Value = 7087.5
$ws.Range("I61").Value = 7087.5
$ws.Range("K61").Value = 21262.5
$ws.Range("M61").Value = -21090.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 3613.0908
$ws.Range("I86").Value = 1492
$ws.Range("K86").Value = 1492
$ws.Range("M86").Value = -369

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H89").Value = 3613.0908
$ws.Range("I89").Value = 1492
$ws.Range("K89").Value = 7460
$ws.Range("M89").Value = -1844

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 10348.5
$ws.Range("I113").Value = 12488.1
$ws.Range("J113").Value = 4999.5
$ws.Range("K113").Value = 12488.1
$ws.Range("L113").Value = 4999.5
$ws.Range("M113").Value = -9234.1
$ws.Range("N113").Value = -11507.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 2296.5
$ws.Range("I135").Value = 2151.6667
$ws.Range("K135").Value = 19365.0003
$ws.Range("M135").Value = -16830.0003

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 45470.92
$ws.Range("I137").Value = 93559.664
$ws.Range("J137").Value = 4252
$ws.Range("K137").Value = 280678.992
$ws.Range("L137").Value = 12756
$ws.Range("M137").Value = -278128.992
$ws.Range("N137").Value = -17856

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 4485.207
$ws.Range("I61").Value = 4373
$ws.Range("J61").Value = 6000
$ws.Range("K61").Value = 4373
$ws.Range("L61").Value = 6000
$ws.Range("M61").Value = -4161
$ws.Range("N61").Value = -6424

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 69833.86
$ws.Range("I74").Value = 69833.86
$ws.Range("K74").Value = 69833.86
$ws.Range("M74").Value = -68959.86

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 69833.86
$ws.Range("I77").Value = 69833.86
$ws.Range("K77").Value = 349169.3
$ws.Range("M77").Value = -344801.3

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 1145.5385
$ws.Range("I97").Value = 952.69696
$ws.Range("K97").Value = 952.69696
$ws.Range("M97").Value = -456.69696

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2970.2273
$ws.Range("I132").Value = 1993.0714
$ws.Range("K132").Value = 5979.2142
$ws.Range("M132").Value = -3449.2142

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 4485.207
$ws.Range("I136").Value = 4373
$ws.Range("J136").Value = 6000
$ws.Range("K136").Value = 13119
$ws.Range("L136").Value = 18000
$ws.Range("M136").Value = -10569
$ws.Range("N136").Value = -23100

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 787.46875
$ws.Range("J80").Value = 494.70587
$ws.Range("L80").Value = 494.70587
$ws.Range("N80").Value = -2490.70587

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H83").Value = 787.46875
$ws.Range("J83").Value = 494.70587
$ws.Range("L83").Value = 2473.52935
$ws.Range("N83").Value = -12457.52935

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 5584.8276
$ws.Range("I134").Value = 5584.8276
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 16754.4828
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -14219.4828
$ws.Range("N134").Value = ""

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H138").Value = 56148.15
$ws.Range("J138").Value = 56148.15
$ws.Range("L138").Value = 56148.15
$ws.Range("N138").Value = -66428.14999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 235153.86
$ws.Range("I31").Value = 265392.62
$ws.Range("J31").Value = 5339.2
$ws.Range("K31").Value = 265392.62
$ws.Range("L31").Value = 5339.2
$ws.Range("M31").Value = -265097.62
$ws.Range("N31").Value = -5929.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 235153.86
$ws.Range("I34").Value = 265392.62
$ws.Range("J34").Value = 5339.2
$ws.Range("K34").Value = 265392.62
$ws.Range("L34").Value = 5339.2
$ws.Range("M34").Value = -265190.62
$ws.Range("N34").Value = -5743.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 3881.1538
$ws.Range("I122").Value = 3880.1667
$ws.Range("J122").Value = 3882
$ws.Range("K122").Value = 11640.5001
$ws.Range("L122").Value = 11646
$ws.Range("M122").Value = -9190.500100000001
$ws.Range("N122").Value = -16546

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 3907.375
$ws.Range("I132").Value = 4050.2432
$ws.Range("K132").Value = 12150.7296
$ws.Range("M132").Value = -9620.729599999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 2222.0715
$ws.Range("I134").Value = 2222.0715
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 6666.2145
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -4131.2145
$ws.Range("N134").Value = ""

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").Value = ""

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 4438
$ws.Range("I80").Value = 3944.5
$ws.Range("J80").Value = 4869.8125
$ws.Range("K80").Value = 3944.5
$ws.Range("L80").Value = 4869.8125
$ws.Range("M80").Value = -2946.5
$ws.Range("N80").Value = -6865.8125

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 4438
$ws.Range("I83").Value = 3944.5
$ws.Range("J83").Value = 4869.8125
$ws.Range("K83").Value = 19722.5
$ws.Range("L83").Value = 24349.0625
$ws.Range("M83").Value = -14730.5
$ws.Range("N83").Value = -34333.0625

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 14277.826
$ws.Range("I126").Value = 15050.632
$ws.Range("J126").Value = 10607
$ws.Range("K126").Value = 45151.896
$ws.Range("L126").Value = 31821
$ws.Range("M126").Value = -42681.896
$ws.Range("N126").Value = -36761

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 25367.623
$ws.Range("I132").Value = 26193.814
$ws.Range("J132").Value = 7604.5
$ws.Range("K132").Value = 78581.442
$ws.Range("L132").Value = 22813.5
$ws.Range("M132").Value = -76051.442
$ws.Range("N132").Value = -27873.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").Value = ""

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 66668376
$ws.Range("I93").Value = 1512.9231
$ws.Range("J93").Value = 500003000
$ws.Range("K93").Value = 1512.9231
$ws.Range("L93").Value = 500003000
$ws.Range("M93").Value = -264.9231
$ws.Range("N93").Value = -500005496

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 5059.625
$ws.Range("I122").Value = 4925.2856
$ws.Range("K122").Value = 14775.8568
$ws.Range("M122").Value = -12325.8568

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 3520.2
$ws.Range("J132").Value = 3145
$ws.Range("L132").Value = 9435
$ws.Range("N132").Value = -14495

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2333.1177
$ws.Range("I122").Value = 2210.9333
$ws.Range("J122").Value = 3249.5
$ws.Range("K122").Value = 6632.7999
$ws.Range("L122").Value = 9748.5
$ws.Range("M122").Value = -4182.7999
$ws.Range("N122").Value = -14648.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 50581.668
$ws.Range("I126").Value = 2009.25
$ws.Range("K126").Value = 6027.75
$ws.Range("M126").Value = -3557.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3751.7646
$ws.Range("J132").Value = 4284
$ws.Range("L132").Value = 12852
$ws.Range("N132").Value = -17912

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1667633.5
$ws.Range("I136").Value = 1667633.5
$ws.Range("K136").Value = 5002900.5
$ws.Range("M136").Value = -5000350.5
